# This script applies the "Deploying to gh-pages" update to the
# StructureDefinition-employee-education-level workbook:
#   - Rebrand from IBM/Alvearie to LinuxForHealth (URL + Publisher)
#   - Bump version 7.0.0 -> 8.0.0
#   - Update the publication Date
#   - Clear the stray/duplicated Constraints text on the "Extension" row
#     of the Elements sheet (it was incorrectly duplicated onto both the
#     Extension row and the Extension.extension row; it should only live
#     on Extension.extension).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Metadata": simple Property/Value pairs
# ---------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-education-level"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------
# Sheet "Elements": per-element definition table
# ---------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 5 is the "Extension.url" element; its Fixed Value (column Q)
# mirrors the StructureDefinition URL and must be updated to match.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-education-level"

# Row 2 is the top-level "Extension" element; its Constraints cell
# (column AI) incorrectly duplicated the ele-1/ext-1 constraint text
# that belongs only on the "Extension.extension" row (row 4). Clear it.
$elements.Range("AI2").Value = ""
